$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (column D) and 1h volume change (column E) values
# D-column values are plain strings that can look numeric (e.g. "1.001"), so we
# force the cell to Text format before assignment and then restore the default
# "Normal" style afterwards so no stray style index is left on the cell (matching
# the source workbook, where these cells carry no explicit style).

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '27.078.93'
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -2.87%  '

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.710.17'
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -3.31%  '

$dCell = $ws.Cells.Item(4, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.16%  '

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '307.60'
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -6.24%  '

$ws.Cells.Item(6, 5).Value = '  -0.05%  '

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.4719'
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +5.52%  '

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.3420'
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -3.73%  '

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '42.13'
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.28%  '

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.07264'
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.19%  '

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.033'
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -6.12%  '

$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.06%  '

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '19.76'
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -5.74%  '

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.836'
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -3.09%  '

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.706.01'
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -3.61%  '

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '6.831'
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -5.62%  '

$ws.Cells.Item(17, 5).Value = '  -4.75%  '

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.00001035'
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -2.41%  '

$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.06357'
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.06%  '

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.04%  '

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '16.43'
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -3.98%  '

$ws.Cells.Item(22, 5).Value = '  -2.97%  '

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '27.102.82'
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -2.98%  '

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '10.83'
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -3.95%  '

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.109'
$dCell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.14%  '

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '156.97'
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -2.93%  '

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '19.45'
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -4.53%  '

$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.904.04'
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -3.55%  '

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.073'
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -4.12%  '

$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '119.12'
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -4.62%  '

$ws.Cells.Item(31, 5).Value = '  -9.03%  '

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.09131'
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.83%  '

$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '3.582'
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -2.52%  '

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.294'
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -6.20%  '

$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.02184'
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -4.61%  '

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.05805'
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -5.89%  '

$dCell = $ws.Cells.Item(37, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '10.99'
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -7.24%  '

$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.1983'
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -5.60%  '

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.9998'
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.14%  '

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '4.724'
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -4.84%  '

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.394'
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.01%  '

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.5864'
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -7.30%  '

$ws.Cells.Item(43, 5).Value = '  -7.77%  '

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '7.458'
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -5.51%  '

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '12.53'
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -5.27%  '

$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.5631'
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -4.19%  '

$ws.Cells.Item(47, 5).Value = '  -5.21%  '

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '117.19'
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -4.42%  '

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.831'
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -6.49%  '

$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.06626'
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -4.03%  '

$ws.Cells.Item(51, 5).Value = '  -4.94%  '
